{"js": "// Apply the \"Added many more features\" edits to the Astro Pug review doc.\n// Each entry is an exact, unique (or intentionally-repeated) old->new text pair\n// taken from the unified diff.\nconst replacements = [\n  {\n    from: \"Play Astro Pug for Free: Chinese Tradition Meets Pugs\",\n    to: \"Play Astro Pug Free: Review of the Chinese-themed Online Slot\",\n  },\n  {\n    from: \"Great graphics and design\",\n    to: \"Beautiful graphics and design\",\n  },\n  {\n    from: \"Exciting free spins feature with all pugs remaining wild\",\n    to: \"Exciting gameplay with free spins and multipliers\",\n  },\n  {\n    from: \"Potential for big wins with the win multiplier\",\n    to: \"Potential for big wins\",\n  },\n  {\n    from: \"A limited number of features compared to other slot games\",\n    to: \"Limited number of unique gameplay features\",\n  },\n  {\n    from: \"Astro Pug may not appeal to players who don't enjoy Chinese themes\",\n    to: \"May not appeal to players who don't like Chinese or pug themes\",\n  },\n  {\n    from: \"Get ready to play Astro Pug for free! With a unique combination of Chinese architecture and pugs, this game features exciting free spins and multipliers.\",\n    to: \"Play Astro Pug for free and enjoy a unique combination of Chinese tradition and adorable pugs. Review of gameplay and features.\",\n  },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Astro Pug review doc.\n$d = $word.ActiveDocument\n\n# Replace every occurrence of $findText with $replaceText. Uses Range.Text (rather than\n# Find.Replacement.Text) so straight apostrophes are not auto-corrected into curly ones,\n# and loops with Wrap = wdFindStop (0) so it keeps finding/replacing every match (the\n# title line occurs twice in this document) without looping forever.\nfunction Replace-AllOccurrences($findText, $replaceText) {\n    $count = 0\n    while ($true) {\n        $rng = $d.Content\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Text = $findText\n        $find.Forward = $true\n        $find.Wrap = 0  # wdFindStop\n        $find.MatchCase = $true\n        $found = $find.Execute()\n        if (-not $found) { break }\n        $rng.Text = $replaceText\n        $count += 1\n        if ($count -gt 25) { break }  # safety net against accidental infinite loop\n    }\n    return $count\n}\n\nReplace-AllOccurrences \"Play Astro Pug for Free: Chinese Tradition Meets Pugs\" \"Play Astro Pug Free: Review of the Chinese-themed Online Slot\" | Out-Null\nReplace-AllOccurrences \"Great graphics and design\" \"Beautiful graphics and design\" | Out-Null\nReplace-AllOccurrences \"Exciting free spins feature with all pugs remaining wild\" \"Exciting gameplay with free spins and multipliers\" | Out-Null\nReplace-AllOccurrences \"Potential for big wins with the win multiplier\" \"Potential for big wins\" | Out-Null\nReplace-AllOccurrences \"A limited number of features compared to other slot games\" \"Limited number of unique gameplay features\" | Out-Null\nReplace-AllOccurrences \"Astro Pug may not appeal to players who don't enjoy Chinese themes\" \"May not appeal to players who don't like Chinese or pug themes\" | Out-Null\nReplace-AllOccurrences \"Get ready to play Astro Pug for free! With a unique combination of Chinese architecture and pugs, this game features exciting free spins and multipliers.\" \"Play Astro Pug for free and enjoy a unique combination of Chinese tradition and adorable pugs. Review of gameplay and features.\" | Out-Null\n"}
